# Update cryptocurrency price/volume data per latest GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.684.31'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '1.571.49'
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").Value = '''210.63'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").Value = '  +6.56%  '
$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D8").Value = '''25.19'
$ws.Range("E8").Value = '  +7.73%  '
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("D11").Value = '''0.0898'
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("D12").Value = '1.796.43'
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("D13").Value = '1.567.55'
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("D14").Value = '28.719.04'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").Value = '''0.518'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '''61.79'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").Value = '''229.99'
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").Value = '''7.31'
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").Value = '0.0₃0687'
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = '''3.95'
$ws.Range("E22").Value = '  +1.43%  '
$ws.Range("D23").Value = '''9.11'
$ws.Range("E23").Value = '  +3.54%  '
$ws.Range("E24").Value = '  +3.94%  '
$ws.Range("D25").Value = '''152.49'
$ws.Range("E26").Value = '  +3.99%  '
$ws.Range("D27").Value = '''14.86'
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("D31").Value = '''1.05'
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("E32").Value = '  +1.30%  '
$ws.Range("D33").Value = '1.408.94'
$ws.Range("D34").Value = '''3.01'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("D36").Value = '''1.48'
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("E37").Value = '  +5.80%  '
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("D40").Value = '''0.519'
$ws.Range("E40").Value = '  +1.84%  '
$ws.Range("D41").Value = '''1.94'
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").Value = '''0.0461'
$ws.Range("E44").Value = '  -0.51%  '
$ws.Range("D45").Value = '''63.62'
$ws.Range("E45").Value = '  +3.49%  '
$ws.Range("D46").Value = '''5.26'
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("D47").Value = '1.707.75'
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("D48").Value = '''84.88'
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").Value = '''0.829'
$ws.Range("E49").Value = '  -8.82%  '
$ws.Range("D50").Value = '''42.19'
$ws.Range("E50").Value = '  +2.07%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.0510'
$ws.Range("E51").Value = '  +0.21%  '
